$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43 (hunk @@ -2760,25 +2760,25 @@)
$ws.Range("H43").Value = 150242
$ws.Range("I43").Value = 5745
$ws.Range("J43").Value = 208040.8
$ws.Range("K43").Value = 5745
$ws.Range("L43").Value = 208040.8
$ws.Range("M43").Value = -5676
$ws.Range("N43").Value = -208178.8

# Row 57 (hunk @@ -3464,19 +3464,22 @@)
$ws.Range("H57").Value = 125750
$ws.Range("J57").Value = 125750
$ws.Range("L57").Value = 377250
$ws.Range("N57").Value = -378248

# Row 132 (hunk @@ -7247,22 +7250,22 @@)
$ws.Range("H132").Value = 3802.4285
$ws.Range("I132").Value = 3802.4285
$ws.Range("K132").Value = 11407.2855
$ws.Range("M132").Value = -8877.2855

# Row 137 (hunk @@ -7495,25 +7498,25 @@)
$ws.Range("H137").Value = 5516
$ws.Range("I137").Value = 1390.4546
$ws.Range("J137").Value = 16861.25
$ws.Range("K137").Value = 4171.3638
$ws.Range("L137").Value = 50583.75
$ws.Range("M137").Value = -1621.3638
$ws.Range("N137").Value = -55683.75

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (hunk @@ -7837,22 +7840,22 @@)
$ws.Range("H2").Value = 1537
$ws.Range("I2").Value = 997.1111
$ws.Range("K2").Value = 997.1111
$ws.Range("M2").Value = -884.1111

# Row 32 (hunk @@ -9310,25 +9313,25 @@)
$ws.Range("H32").Value = 10640729
$ws.Range("J32").Value = 11753
$ws.Range("L32").Value = 11753
$ws.Range("N32").Value = -12327

# Row 61 (hunk @@ -10707,22 +10710,22 @@)
$ws.Range("H61").Value = 20051660
$ws.Range("I61").Value = 23820228
$ws.Range("K61").Value = 23820228
$ws.Range("M61").Value = -23820016

# Row 74 (hunk @@ -11335,25 +11338,25 @@)
$ws.Range("H74").Value = 8072182
$ws.Range("I74").Value = 11906542
$ws.Range("J74").Value = 20026.8
$ws.Range("K74").Value = 11906542
$ws.Range("L74").Value = 20026.8
$ws.Range("M74").Value = -11905668
$ws.Range("N74").Value = -21774.8

# Row 77 (hunk @@ -11482,25 +11485,25 @@)
$ws.Range("H77").Value = 8072182
$ws.Range("I77").Value = 11906542
$ws.Range("J77").Value = 20026.8
$ws.Range("K77").Value = 59532710
$ws.Range("L77").Value = 100134
$ws.Range("M77").Value = -59528342
$ws.Range("N77").Value = -108870

# Row 110 (hunk @@ -13099,22 +13102,22 @@)
$ws.Range("H110").Value = 1972
$ws.Range("I110").Value = 1972
$ws.Range("K110").Value = 1972
$ws.Range("M110").Value = 73

# Row 116 (hunk @@ -13390,22 +13393,22 @@)
$ws.Range("H116").Value = 1537
$ws.Range("I116").Value = 997.1111
$ws.Range("K116").Value = 997.1111
$ws.Range("M116").Value = 1296.8889

# Row 136 (hunk @@ -14373,22 +14376,22 @@)
$ws.Range("H136").Value = 20051660
$ws.Range("I136").Value = 23820228
$ws.Range("K136").Value = 71460684
$ws.Range("M136").Value = -71458134

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (hunk @@ -14819,22 +14822,22 @@)
$ws.Range("H3").Value = 1537
$ws.Range("I3").Value = 997.1111
$ws.Range("K3").Value = 997.1111
$ws.Range("M3").Value = -883.1111

# Row 80 (hunk @@ -18517,25 +18520,25 @@)
$ws.Range("H80").Value = 3004.6
$ws.Range("J80").Value = 1837.25
$ws.Range("L80").Value = 1837.25
$ws.Range("N80").Value = -3833.25

# Row 83 (hunk @@ -18670,25 +18673,25 @@)
$ws.Range("H83").Value = 3004.6
$ws.Range("J83").Value = 1837.25
$ws.Range("L83").Value = 9186.25
$ws.Range("N83").Value = -19170.25

# Row 105 (hunk @@ -19781,25 +19784,25 @@)
$ws.Range("H105").Value = 1817
$ws.Range("I105").Value = 1018.9
$ws.Range("J105").Value = 2957.1428
$ws.Range("K105").Value = 1018.9
$ws.Range("L105").Value = 2957.1428
$ws.Range("M105").Value = 728.1
$ws.Range("N105").Value = -6451.1428

# Row 134 (hunk @@ -21208,22 +21211,22 @@)
$ws.Range("H134").Value = 82724.234
$ws.Range("I134").Value = 1538.7
$ws.Range("K134").Value = 4616.1
$ws.Range("M134").Value = -2081.1

$ws = $wb.Worksheets.Item("CRP")
# Row 8 (hunk @@ -21994,25 +21997,19 @@)
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()

# Row 31 (hunk @@ -23121,25 +23118,25 @@)
$ws.Range("H31").Value = 358957.97
$ws.Range("I31").Value = 4765.6743
$ws.Range("J31").Value = 1021143.5
$ws.Range("K31").Value = 4765.6743
$ws.Range("L31").Value = 1021143.5
$ws.Range("M31").Value = -4470.6743
$ws.Range("N31").Value = -1021733.5

# Row 34 (hunk @@ -23271,25 +23268,25 @@)
$ws.Range("H34").Value = 358957.97
$ws.Range("I34").Value = 4765.6743
$ws.Range("J34").Value = 1021143.5
$ws.Range("K34").Value = 4765.6743
$ws.Range("L34").Value = 1021143.5
$ws.Range("M34").Value = -4563.6743
$ws.Range("N34").Value = -1021547.5

# Row 51 (hunk @@ -24098,22 +24095,22 @@)
$ws.Range("H51").Value = 24999.666
$ws.Range("I51").Value = 24999.666
$ws.Range("K51").Value = 24999.666
$ws.Range("M51").Value = -24263.666

# Row 61 (hunk @@ -24582,22 +24579,22 @@)
$ws.Range("H61").Value = 24999.666
$ws.Range("I61").Value = 24999.666
$ws.Range("K61").Value = 24999.666
$ws.Range("M61").Value = -24651.666

# Row 99 (hunk @@ -26426,25 +26423,25 @@)
$ws.Range("H99").Value = 3508.0386
$ws.Range("J99").Value = 4062.125
$ws.Range("L99").Value = 4062.125
$ws.Range("N99").Value = -7058.125

# Row 107 (hunk @@ -26827,25 +26824,25 @@)
$ws.Range("H107").Value = 1609
$ws.Range("J107").Value = 3174.75
$ws.Range("L107").Value = 3174.75
$ws.Range("N107").Value = -7014.75

# Row 122 (hunk @@ -27565,22 +27562,22 @@)
$ws.Range("H122").Value = 2226.8235
$ws.Range("I122").Value = 2226.8235
$ws.Range("K122").Value = 6680.470499999999
$ws.Range("M122").Value = -4230.470499999999

# Row 126 (hunk @@ -27761,25 +27758,25 @@)
$ws.Range("H126").Value = 3508.0386
$ws.Range("J126").Value = 4062.125
$ws.Range("L126").Value = 12186.375
$ws.Range("N126").Value = -17126.375

# Row 134 (hunk @@ -28159,22 +28156,22 @@)
$ws.Range("H134").Value = 589832.0600000001
$ws.Range("I134").Value = 589832.0600000001
$ws.Range("K134").Value = 1769496.18
$ws.Range("M134").Value = -1766961.18

$ws = $wb.Worksheets.Item("CUL")
# Row 38 (hunk @@ -30469,25 +30466,25 @@)
$ws.Range("H38").Value = 166.06667
$ws.Range("I38").Value = 297.42856
$ws.Range("J38").Value = 51.125
$ws.Range("K38").Value = 892.28568
$ws.Range("L38").Value = 153.375
$ws.Range("M38").Value = -545.28568
$ws.Range("N38").Value = -847.375

# Row 64 (hunk @@ -31788,25 +31785,25 @@)
$ws.Range("H64").Value = 3866.3333
$ws.Range("I64").Value = 3399
$ws.Range("J64").Value = 4100
$ws.Range("K64").Value = 10197
$ws.Range("L64").Value = 12300
$ws.Range("M64").Value = -9927
$ws.Range("N64").Value = -12840

# Row 67 (hunk @@ -31932,25 +31929,25 @@)
$ws.Range("H67").Value = 3866.3333
$ws.Range("I67").Value = 3399
$ws.Range("J67").Value = 4100
$ws.Range("K67").Value = 10197
$ws.Range("L67").Value = 12300
$ws.Range("M67").Value = -9261
$ws.Range("N67").Value = -14172

# Row 86 (hunk @@ -32884,25 +32881,25 @@)
$ws.Range("H86").Value = 566.1667
$ws.Range("I86").Value = 550
$ws.Range("J86").Value = 574.25
$ws.Range("K86").Value = 1650
$ws.Range("L86").Value = 1722.75
$ws.Range("M86").Value = -464
$ws.Range("N86").Value = -4094.75

# Row 89 (hunk @@ -33037,25 +33034,25 @@)
$ws.Range("H89").Value = 566.1667
$ws.Range("I89").Value = 550
$ws.Range("J89").Value = 574.25
$ws.Range("K89").Value = 4950
$ws.Range("L89").Value = 5168.25
$ws.Range("M89").Value = 978
$ws.Range("N89").Value = -17024.25

# Row 123 (hunk @@ -34754,25 +34751,25 @@)
$ws.Range("H123").Value = 5333
$ws.Range("I123").Value = 7000
$ws.Range("J123").Value = 4999.6
$ws.Range("K123").Value = 21000
$ws.Range("L123").Value = 14998.8
$ws.Range("M123").Value = -18550
$ws.Range("N123").Value = -19898.8

$ws = $wb.Worksheets.Item("GSM")
# Row 119 (hunk @@ -41524,22 +41521,22 @@)
$ws.Range("H119").Value = 112000
$ws.Range("J119").Value = 112000
$ws.Range("L119").Value = 112000
$ws.Range("N119").Value = -121676

# Row 132 (hunk @@ -42158,22 +42155,22 @@)
$ws.Range("H132").Value = 166672480
$ws.Range("I132").Value = 250007730
$ws.Range("K132").Value = 750023190
$ws.Range("M132").Value = -750020660

$ws = $wb.Worksheets.Item("LTW")
# Row 122 (hunk @@ -48589,25 +48586,25 @@)
$ws.Range("H122").Value = 5563.6
$ws.Range("I122").Value = 5245.3438
$ws.Range("J122").Value = 6836.625
$ws.Range("K122").Value = 15736.0314
$ws.Range("L122").Value = 20509.875
$ws.Range("M122").Value = -13286.0314
$ws.Range("N122").Value = -25409.875

# Row 123 (hunk @@ -48641,22 +48638,22 @@)
$ws.Range("H123").Value = 161000
$ws.Range("J123").Value = 161000
$ws.Range("L123").Value = 161000
$ws.Range("N123").Value = -170800

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (hunk @@ -55561,22 +55558,22 @@)
$ws.Range("H122").Value = 13807.321
$ws.Range("I122").Value = 9548.315000000001
$ws.Range("K122").Value = 28644.945
$ws.Range("M122").Value = -26194.945

# Row 126 (hunk @@ -55760,25 +55757,25 @@)
$ws.Range("H126").Value = 5395.033
$ws.Range("I126").Value = 4887.185
$ws.Range("J126").Value = 9965.666999999999
$ws.Range("K126").Value = 14661.555
$ws.Range("L126").Value = 29897.001
$ws.Range("M126").Value = -12191.555
$ws.Range("N126").Value = -34837.001

# Row 132 (hunk @@ -56057,22 +56054,22 @@)
$ws.Range("H132").Value = 1715.1538
$ws.Range("I132").Value = 1688.258
$ws.Range("K132").Value = 5064.774
$ws.Range("M132").Value = -2534.774
